# Regenerate save_data to use K instead of Strike# in column G (header "K").
# Update the G column values for rows 2-15 with the newly computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 6
    4  = 7
    5  = 4
    6  = 1
    7  = 0
    8  = 5
    9  = 6
    10 = 4
    11 = 6
    12 = 1
    13 = 6
    14 = 0
    15 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
